$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab7")

$ws.Range("B3").Value = 6250.6
$ws.Range("C3").Value = 3237.44
$ws.Range("B4").Value = 1250.3
$ws.Range("C4").Value = 870.06
$ws.Range("B5").Value = 892.875
$ws.Range("C5").Value = 627.61
$ws.Range("B6").Value = 694.5
$ws.Range("C6").Value = 621.82000000000005

$ws.Range("B11").Value = 546.79999999999995
$ws.Range("C11").Value = 280.14
$ws.Range("B12").Value = 273.39999999999998
$ws.Range("C12").Value = 212.45
$ws.Range("B13").Value = 182.32
$ws.Range("C13").Value = 182.99
$ws.Range("B14").Value = 136.71
$ws.Range("C14").Value = 185.24
